$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "first"
$ws2.Name = "second"

# --- Sheet "first": row 2 (A2:C2) switches from style 2 to style 1 ---
# A1 already carries style 1; copy its format (not value) onto A2:C2.
$ws1.Range("A1").Copy()
$ws1.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet "second" ---
# Keep a style-2 source cell before we touch A1..A6 (A2 keeps style 2 unchanged).
# Apply style 1 (copied from B1, which is already style 1) to A1,A3,A4,A5,A6.
$ws2.Range("B1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A3").PasteSpecial(-4122)
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Range("A5").PasteSpecial(-4122)
$ws2.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# B4 switches from style 1 to style 2; copy format from A2 (still style 2).
$ws2.Range("A2").Copy()
$ws2.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update values on sheet "second" ---
$ws2.Range("B1").Value = 20.0
$ws2.Range("B2").Value = 10.0
$ws2.Range("B3").Value = 30.0
$ws2.Range("B4").Value = -30.0
$ws2.Range("B5").Value = 35.0
$ws2.Range("B6").Value = 40.0
$ws2.Range("B7").Value = 100.0
$ws2.Range("B8").Value = 101.0
$ws2.Range("B9").Value = 101.0
$ws2.Range("B10").Value = 101.0
$ws2.Range("B11").Value = 101.0

# --- New rows 12-14, with style 1 (copied from A7:B7) ---
$ws2.Range("A12").Value = 1.0
$ws2.Range("B12").Value = 101.0
$ws2.Range("A13").Value = 1.0
$ws2.Range("B13").Value = -300.0
$ws2.Range("A14").Value = 1.0
$ws2.Range("B14").Value = 400.0

$ws2.Range("A7:B7").Copy()
$ws2.Range("A12:B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
